$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.02456200430219724
$ws.Range("D2").Value = 0.024470825559753
$ws.Range("E2").Value = 0.4248831183909374
$ws.Range("F2").Value = 0.3464024739650284
$ws.Range("G2").Value = 0.2005477842955941
$ws.Range("H2").Value = 0.379178972343901
$ws.Range("I2").Value = 0.2515818745413299
$ws.Range("K2").Value = 1.784400975634981
$ws.Range("O2").Value = 1.068191922977817

$ws.Range("C3").Value = 0.02143509417579992
$ws.Range("D3").Value = 0.0214550273307097
$ws.Range("E3").Value = 0.3706405723224719
$ws.Range("F3").Value = 0.3480952699450341
$ws.Range("G3").Value = 0.2030354731723278
$ws.Range("H3").Value = 0.3850595669738652
$ws.Range("I3").Value = 0.2535577513752401
$ws.Range("K3").Value = 1.557906378721498
$ws.Range("O3").Value = 1.085486744642807

$ws.Range("C4").Value = 0.01950949172042016
$ws.Range("D4").Value = 0.01959489452532637
$ws.Range("E4").Value = 0.3374208952504603
$ws.Range("F4").Value = 0.3495513776082149
$ws.Range("G4").Value = 0.2049025081595843
$ws.Range("H4").Value = 0.3889781425419869
$ws.Range("I4").Value = 0.25508141856427
$ws.Range("K4").Value = 1.418251538986794
$ws.Range("O4").Value = 1.097460660044987

$ws.Range("C5").Value = 0.01872341666005894
$ws.Range("D5").Value = 0.0188348175072548
$ws.Range("E5").Value = 0.3239034082975962
$ws.Range("F5").Value = 0.3502490677200925
$ws.Range("G5").Value = 0.2057481309850999
$ws.Range("H5").Value = 0.3906522205825596
$ws.Range("I5").Value = 0.2557800282795633
$ws.Range("K5").Value = 1.361197714077605
$ws.Range("O5").Value = 1.102679135351053

$ws.Range("C6").Value = 0.0185928078027402
$ws.Range("D6").Value = 0.01870848477923914
$ws.Range("E6").Value = 0.3216599866585739
$ws.Range("F6").Value = 0.3503712050539818
$ws.Range("G6").Value = 0.2058936495703634
$ws.Range("H6").Value = 0.3909348587326207
$ws.Range("I6").Value = 0.2559007141807328
$ws.Range("K6").Value = 1.351715433141521
$ws.Range("O6").Value = 1.103566085333895

$ws.Range("C7").Value = 0.01949889594731502
$ws.Range("D7").Value = 0.01958465211271232
$ws.Range("E7").Value = 0.3372385157348674
$ws.Range("F7").Value = 0.349560365166802
$ws.Range("G7").Value = 0.2049135700071218
$ws.Range("H7").Value = 0.3890004073086004
$ws.Range("I7").Value = 0.255090526131589
$ws.Range("K7").Value = 1.417482666108356
$ws.Range("O7").Value = 1.097529667826521

$ws.Range("C8").Value = 0.02348504994654377
$ws.Range("D8").Value = 0.02343276064935651
$ws.Range("E8").Value = 0.4061615152260032
$ws.Range("F8").Value = 0.346899373121758
$ws.Range("G8").Value = 0.2013346979065673
$ws.Range("H8").Value = 0.3811425904529813
$ws.Range("I8").Value = 0.252198497421027
$ws.Range("K8").Value = 1.706429671766102
$ws.Range("O8").Value = 1.073872977814048

$ws.Range("C9").Value = 0.03125517362950347
$ws.Range("D9").Value = 0.03090983621817145
$ws.Range("E9").Value = 0.5420800674901614
$ws.Range("F9").Value = 0.3450094610152519
$ws.Range("G9").Value = 0.1970380628082467
$ws.Range("H9").Value = 0.3681846200660317
$ws.Range("I9").Value = 0.2490074440067858
$ws.Range("K9").Value = 2.268252045987879
$ws.Range("O9").Value = 1.038307482663711

$ws.Range("C10").Value = 0.03693362884746421
$ws.Range("D10").Value = 0.03635869126316038
$ws.Range("E10").Value = 0.6425268476959758
$ws.Range("F10").Value = 0.3456806227796108
$ws.Range("G10").Value = 0.1955779129907924
$ws.Range("H10").Value = 0.3601701640416124
$ws.Range("I10").Value = 0.2481984410142886
$ws.Range("K10").Value = 2.677935760104219
$ws.Range("O10").Value = 1.018880683339091

$ws.Range("C11").Value = 0.0395100001719868
$ws.Range("D11").Value = 0.03882736665640607
$ws.Range("E11").Value = 0.6883776182536252
$ws.Range("F11").Value = 0.3464398041583152
$ws.Range("G11").Value = 0.1952900962870885
$ws.Range("H11").Value = 0.3568536073079471
$ws.Range("I11").Value = 0.2481688688149504
$ws.Range("K11").Value = 2.863610638929913
$ws.Range("O11").Value = 1.011520760251813

$ws.Range("C12").Value = 0.04048458906603969
$ws.Range("D12").Value = 0.03976069612052413
$ws.Range("E12").Value = 0.7057647355573522
$ws.Range("F12").Value = 0.3467930748156789
$ws.Range("G12").Value = 0.1952358855262304
$ws.Range("H12").Value = 0.3556452687326868
$ws.Range("I12").Value = 0.2482067422034646
$ws.Range("K12").Value = 2.933817967073026
$ws.Range("O12").Value = 1.008948021356048

$ws.Range("C13").Value = 0.04027474039560275
$ws.Range("D13").Value = 0.03955975479213691
$ws.Range("E13").Value = 0.7020189928327625
$ws.Range("F13").Value = 0.346714057582254
$ws.Range("G13").Value = 0.1952451140971974
$ws.Range("H13").Value = 0.3559033868554167
$ws.Range("I13").Value = 0.2481963965953611
$ws.Range("K13").Value = 2.918702241447079
$ws.Range("O13").Value = 1.009492545982994

$ws.Range("C14").Value = 0.03959020106158562
$ws.Range("D14").Value = 0.03890418272979446
$ws.Range("E14").Value = 0.6898075666943413
$ws.Range("F14").Value = 0.3464675465426339
$ws.Range("G14").Value = 0.1952845347494048
$ws.Range("H14").Value = 0.3567532419286366
$ws.Range("I14").Value = 0.2481709989022249
$ws.Range("K14").Value = 2.869388740414081
$ws.Range("O14").Value = 1.011304793891654

$ws.Range("C15").Value = 0.03917076569707945
$ws.Range("D15").Value = 0.03850242820305994
$ws.Range("E15").Value = 0.6823309516142899
$ws.Range("F15").Value = 0.3463251339864257
$ws.Range("G15").Value = 0.1953158342275287
$ws.Range("H15").Value = 0.3572800042542781
$ws.Range("I15").Value = 0.2481618446327332
$ws.Range("K15").Value = 2.839169167190903
$ws.Range("O15").Value = 1.012442811760295

$ws.Range("C16").Value = 0.03676511597967647
$ws.Range("D16").Value = 0.03619715032350257
$ws.Range("E16").Value = 0.6395337091026505
$ws.Range("F16").Value = 0.345640185895256
$ws.Range("G16").Value = 0.1956043525814479
$ws.Range("H16").Value = 0.360393553369704
$ws.Range("I16").Value = 0.2482072185533681
$ws.Range("K16").Value = 2.665787181759015
$ws.Range("O16").Value = 1.019391559555785

$ws.Range("C17").Value = 0.03528755309433507
$ws.Range("D17").Value = 0.03478032276693455
$ws.Range("E17").Value = 0.6133205123348944
$ws.Range("F17").Value = 0.3453366086508112
$ws.Range("G17").Value = 0.1958782213386669
$ws.Range("H17").Value = 0.3623881043109805
$ws.Range("I17").Value = 0.2483220257202241
$ws.Range("K17").Value = 2.559242652070395
$ws.Range("O17").Value = 1.024034113582161

$ws.Range("C18").Value = 0.03443706225732512
$ws.Range("D18").Value = 0.03396445999031528
$ws.Range("E18").Value = 0.5982580435320699
$ws.Range("F18").Value = 0.345204688113192
$ws.Range("G18").Value = 0.1960711342506798
$ws.Range("H18").Value = 0.3635662954074377
$ws.Range("I18").Value = 0.2484198919521781
$ws.Range("K18").Value = 2.49789613968926
$ws.Range("O18").Value = 1.026843339447936

$ws.Range("C19").Value = 0.03414899305057872
$ws.Range("D19").Value = 0.03368806312460038
$ws.Range("E19").Value = 0.5931606138200607
$ws.Range("F19").Value = 0.3451673381233462
$ws.Range("G19").Value = 0.1961425098004668
$ws.Range("H19").Value = 0.3639705244268612
$ws.Range("I19").Value = 0.2484584826841036
$ws.Range("K19").Value = 2.477114243059361
$ws.Range("O19").Value = 1.027818303681116

$ws.Range("C20").Value = 0.03544490835584213
$ws.Range("D20").Value = 0.03493124431884098
$ws.Range("E20").Value = 0.6161094204397415
$ws.Range("F20").Value = 0.3453645030393204
$ws.Range("G20").Value = 0.1958453999999392
$ws.Range("H20").Value = 0.3621725732240861
$ws.Range("I20").Value = 0.2483065068638375
$ws.Range("K20").Value = 2.570591248169706
$ws.Range("O20").Value = 1.023525511385571

$ws.Range("C21").Value = 0.03979129513538737
$ws.Range("D21").Value = 0.03909678144599127
$ws.Range("E21").Value = 0.6933936808570422
$ws.Range("F21").Value = 0.3465381632839524
$ws.Range("G21").Value = 0.1952714640682132
$ws.Range("H21").Value = 0.3565023260976901
$ws.Range("I21").Value = 0.2481771238154309
$ws.Range("K21").Value = 2.883876158867679
$ws.Range("O21").Value = 1.010766661837664

$ws.Range("C22").Value = 0.0426259005443228
$ws.Range("D22").Value = 0.04181040454751894
$ws.Range("E22").Value = 0.744046653726997
$ws.Range("F22").Value = 0.3476889047192557
$ws.Range("G22").Value = 0.1952159606425994
$ws.Range("H22").Value = 0.353073866529023
$ws.Range("I22").Value = 0.2483787789401575
$ws.Range("K22").Value = 3.088019512724827
$ws.Range("O22").Value = 1.003678004296972

$ws.Range("C23").Value = 0.04111358493980788
$ws.Range("D23").Value = 0.04036291768456124
$ws.Range("E23").Value = 0.7169984875258564
$ws.Range("F23").Value = 0.3470394535671986
$ws.Range("G23").Value = 0.1952161256323635
$ws.Range("H23").Value = 0.3548782465859333
$ws.Range("I23").Value = 0.2482448300341034
$ws.Range("K23").Value = 2.979121214591146
$ws.Range("O23").Value = 1.00734636843525

$ws.Range("C24").Value = 0.03537377121642749
$ws.Range("D24").Value = 0.03486301676061032
$ws.Range("E24").Value = 0.614848530818179
$ws.Range("F24").Value = 0.3453517592959798
$ws.Range("G24").Value = 0.1958601281130541
$ws.Range("H24").Value = 0.3622699167295522
$ws.Range("I24").Value = 0.2483134237193809
$ws.Range("K24").Value = 2.565460836138357
$ws.Range("O24").Value = 1.023755013979425

$ws.Range("C25").Value = 0.02915833581940319
$ws.Range("D25").Value = 0.02889473869366554
$ws.Range("E25").Value = 0.5052163354439756
$ws.Range("F25").Value = 0.3451612423365304
$ws.Range("G25").Value = 0.197905182563531
$ws.Range("H25").Value = 0.3714264947906116
$ws.Range("I25").Value = 0.2496027574982698
$ws.Range("K25").Value = 2.116795026417037
$ws.Range("O25").Value = 1.046759242104642

Write-Output "updated 24 rows (216 cells)"